$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking price strings
# (e.g. "148.39") are stored as literal text, matching the inlineStr source.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range('D2').Value2 = '63.065.49'
$ws.Range('E2').Value2 = '  +3.05%  '
$ws.Range('D3').Value2 = '2.953.08'
$ws.Range('E3').Value2 = '  +1.07%  '
$ws.Range('E4').Value2 = '  +0.11%  '
$ws.Range('D5').Value2 = '594.75'
$ws.Range('E5').Value2 = '  -0.37%  '
$ws.Range('D6').Value2 = '148.39'
$ws.Range('E6').Value2 = '  +2.62%  '
$ws.Range('E7').Value2 = '  +0.00%  '
$ws.Range('D8').Value2 = '2.951.90'
$ws.Range('E8').Value2 = '  +1.11%  '
$ws.Range('E9').Value2 = '  +1.45%  '
$ws.Range('D10').Value2 = '7.19'
$ws.Range('E10').Value2 = '  +3.86%  '
$ws.Range('E11').Value2 = '  +6.23%  '
$ws.Range('E12').Value2 = '  +0.86%  '
$ws.Range('D13').Value2 = '0.0000234'
$ws.Range('E13').Value2 = '  +4.88%  '
$ws.Range('D14').Value2 = '32.80'
$ws.Range('E14').Value2 = '  -1.74%  '
$ws.Range('D16').Value2 = '3.445.33'
$ws.Range('E16').Value2 = '  +1.13%  '
$ws.Range('D17').Value2 = '63.063.43'
$ws.Range('E17').Value2 = '  +3.03%  '
$ws.Range('D18').Value2 = '6.69'
$ws.Range('E18').Value2 = '  +0.45%  '
$ws.Range('D19').Value2 = '2.956.43'
$ws.Range('E19').Value2 = '  +1.07%  '
$ws.Range('D20').Value2 = '442.10'
$ws.Range('E20').Value2 = '  +2.53%  '
$ws.Range('D21').Value2 = '13.51'
$ws.Range('E21').Value2 = '  +0.14%  '
$ws.Range('D22').Value2 = '0.669'
$ws.Range('E22').Value2 = '  -0.66%  '
$ws.Range('D23').Value2 = '7.03'
$ws.Range('E23').Value2 = '  -0.53%  '
$ws.Range('D24').Value2 = '11.26'
$ws.Range('E24').Value2 = '  +3.36%  '
$ws.Range('D25').Value2 = '80.87'
$ws.Range('E25').Value2 = '  -1.00%  '
$ws.Range('D26').Value2 = '2.14'
$ws.Range('E26').Value2 = '  -1.83%  '
$ws.Range('D27').Value2 = '11.82'
$ws.Range('E27').Value2 = '  +0.71%  '
$ws.Range('E28').Value2 = '  +0.05%  '
$ws.Range('B29').Value2 = 'ImmutableX'
$ws.Range('C29').Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').Value2 = '2.24'
$ws.Range('E29').Value2 = '  +1.77%  '
$ws.Range('B30').Value2 = 'NEARProtocol'
$ws.Range('C30').Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').Value2 = '7.30'
$ws.Range('E30').Value2 = '  +5.98%  '
$ws.Range('E31').Value2 = '  +0.65%  '
$ws.Range('D32').Value2 = '0.0000101'
$ws.Range('E32').Value2 = '  +15.39%  '
$ws.Range('D33').Value2 = '26.45'
$ws.Range('E33').Value2 = '  -0.48%  '
$ws.Range('E34').Value2 = '  -0.57%  '
$ws.Range('E35').Value2 = '  -0.13%  '
$ws.Range('E36').Value2 = '  -1.60%  '
$ws.Range('E37').Value2 = '  -0.23%  '
$ws.Range('D38').Value2 = '3.08'
$ws.Range('E38').Value2 = '  +3.42%  '
$ws.Range('B39').Value2 = 'OKB'
$ws.Range('C39').Value2 = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').Value2 = '49.77'
$ws.Range('E39').Value2 = '  -0.11%  '
$ws.Range('B40').Value2 = 'Stacks'
$ws.Range('C40').Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value2 = '2.04'
$ws.Range('E40').Value2 = '  +2.39%  '
$ws.Range('D41').Value2 = '8.50'
$ws.Range('E41').Value2 = '  -0.55%  '
$ws.Range('D42').Value2 = '0.118'
$ws.Range('E42').Value2 = '  -3.76%  '
$ws.Range('D43').Value2 = '0.281'
$ws.Range('E43').Value2 = '  +0.33%  '
$ws.Range('D44').Value2 = '39.38'
$ws.Range('E44').Value2 = '  -7.34%  '
$ws.Range('D45').Value2 = '2.708.77'
$ws.Range('E45').Value2 = '  +0.62%  '
$ws.Range('D46').Value2 = '135.66'
$ws.Range('E46').Value2 = '  +1.59%  '
$ws.Range('D47').Value2 = '0.0339'
$ws.Range('E47').Value2 = '  -1.45%  '
$ws.Range('D48').Value2 = '361.40'
$ws.Range('E48').Value2 = '  -0.77%  '
$ws.Range('E49').Value2 = '  +0.01%  '
$ws.Range('D50').Value2 = '0.105'
$ws.Range('E50').Value2 = '  -0.19%  '
$ws.Range('E51').Value2 = '  -3.18%  '

# Restore the original (default) cell style on column D so no stray
# number-format attribute is left behind on any cell.
$dRange.Style = "Normal"

